$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "legia"
$ws.Range("A3").Value = "to"
$ws.Range("A4").Value = "chuje"
$ws.Range("A5").Value = "a"
$ws.Range("A6").Value = "lech"
$ws.Range("A7").Value = "mistrz"
$ws.Range("A8").Value = "polski"

$ws.Range("K17").Select() | Out-Null
